$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.306.89'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '3.159.05'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''593.21'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '''146.85'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D8').Value = '3.147.77'
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').Value = '''0.165'
$ws.Range('E10').Value = '  +3.78%  '
$ws.Range('D11').Value = '''5.90'
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('D12').Value = '''0.458'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '''37.43'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '3.681.58'
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '''7.28'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '64.102.91'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '3.155.79'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').Value = '''468.89'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').Value = '''14.38'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('D22').Value = '''0.735'
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').Value = '''7.61'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').Value = '''2.38'
$ws.Range('E24').Value = '  +11.98%  '
$ws.Range('D25').Value = '''13.17'
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').Value = '''81.27'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '''9.82'
$ws.Range('E28').Value = '  +9.68%  '
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').Value = '''2.23'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').Value = '''7.32'
$ws.Range('E31').Value = '  +7.60%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = '''27.84'
$ws.Range('E33').Value = '  +4.10%  '
$ws.Range('E34').Value = '  +6.07%  '
$ws.Range('D35').Value = '0.0₃0873'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('D37').Value = '''6.17'
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '''3.27'
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('D40').Value = '''467.48'
$ws.Range('E40').Value = '  +6.54%  '
$ws.Range('D41').Value = '''9.38'
$ws.Range('E41').Value = '  +7.46%  '
$ws.Range('D42').Value = '''51.37'
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('E43').Value = '  +9.68%  '
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('D45').Value = '2.912.97'
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('D46').Value = '''39.91'
$ws.Range('E46').Value = '  +11.26%  '
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').Value = '''133.20'
$ws.Range('E48').Value = '  +7.70%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('E51').Value = '  +4.21%  '
